# Apply cryptos list update (coin price/volume refresh, with a couple of rows re-sorted)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.683.67'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '1.596.24'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.39'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.513'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.43'
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0840'
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").Value = '1.820.05'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.600.92'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '26.645.78'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = '0.0₃0750'
$ws.Range("E18").Value = '  +3.02%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.22'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.98'
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  -1.70%  '
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.87'
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.32'
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0517'
$ws.Range("E30").Value = '  +2.34%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").Value = '1.284.95'
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("E35").Value = '  -6.80%  '
$ws.Range("E36").Value = '  -0.51%  '
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.826'
$ws.Range("E39").Value = '  -1.97%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.05'
$ws.Range("E40").Value = '  +16.02%  '
$ws.Range("E41").Value = '  +0.30%  '
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.86'
$ws.Range("E44").Value = '  -1.06%  '
$ws.Range("D45").Value = '1.733.02'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.76'
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.57'
$ws.Range("E47").Value = '  -2.47%  '
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0509'
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.32'
$ws.Range("E51").Value = '  -2.18%  '
